$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets ---
$wb.Worksheets.Item(2).Name = "Include ValueSet #0"
$wb.Worksheets.Item(3).Name = "Include #1"

# --- 2. Update Metadata sheet ---
$ws1 = $wb.Worksheets.Item(1)

# Update the "Date" value (row 8, column B)
$ws1.Cells.Item(8, 2).Value = "2024-09-17T19:55:11+00:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10) and before "Description" (row 11).
# Shift rows 11..14 down to 12..15, carrying both formatting and values, working from the
# bottom up so that source data is not overwritten before it is copied.
for ($r = 14; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws1.Range("A" + $r + ":B" + $r).Copy()
    $ws1.Range("A" + $destRow + ":B" + $destRow).PasteSpecial(-4122)
    $ws1.Cells.Item($destRow, 1).Value = $ws1.Cells.Item($r, 1).Value2
    $ws1.Cells.Item($destRow, 2).Value = $ws1.Cells.Item($r, 2).Value2
}

# Populate the newly freed row 11 with the Jurisdiction property (empty value).
$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = ""
